$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comforter-cda")
$tbl = $ws.ListObjects.Item("comforter_cda_table")

# Grow the table by one row (this keeps the table ref / autoFilter / sheet
# dimension in sync, mirroring what Excel does when a new row of data is
# entered at the bottom of a table).
$tbl.ListRows.Add() | Out-Null

# Row 126 had been a "blank" data row (only the calculated columns were
# present). Fill in the Date / Start Time / End Time that make it a real
# daily power record.
$ws.Range("A126").Value = 43448
$ws.Range("B126").Value = 0.39999999999999997
$ws.Range("C126").Value = 0.69513888888888886

# Row 127: a brand-new daily power record with only a start time so far
# (no end time yet), matching the commit "Add daily power records".
$ws.Range("A127").Value = 43449
$ws.Range("B127").Value = 0.75624999999999998

$ws.Range("D127").Formula = "=(C127-B127)* 1440"
$ws.Range("E127").Formula = "=IF(C127>B127, (C127-B127)*1440, (B127-C127)*1440)"
$ws.Range("F127").Formula = "=ABS((C127-B127)*1440)"

# Reflect the author's final scroll position / selection (the view had
# scrolled down and the new C127 cell is selected after entering data).
$ws.Activate() | Out-Null
$ws.Range("A119").Select() | Out-Null
$ws.Range("C127").Select() | Out-Null
